# Updated test and jenkinsfile
# Update the last row of the messageData fixture: the "subject" test case
# moves from a 12-character subject to a >100-character subject, and the
# related description / expectedMessage columns are filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").Value = "Subject must be between 5 and 100 characters"
$ws.Range("A11").Value = "Testing api with a subject that has more than 100 characters"
$ws.Range("E11").Value = "Testing the api with a very long subject written in the text box field to see if it can fit more that 30 characters"
$ws.Range("F11").Value = "Testing the api for sending message description."

# Move the visible selection (sheetView scrolled one column right, active cell
# moved from G14 to G12).
$ws.Range("G12").Select()
